$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '67.640.74'
$ws.Range("E2").Value = '  -1.23%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.783.26'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.40'
$ws.Range("E5").Value = '  +0.42%  '

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.99'
$ws.Range("E6").Value = '  +0.12%  '

# Row 7: LidoStakedEther
$ws.Range("D7").Value = '3.781.20'
$ws.Range("E7").Value = '  +0.72%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.01%  '

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'

# Row 11: Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("E11").Value = '  -1.82%  '

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  +0.09%  '

# Row 13: ShibaInu
$ws.Range("E13").Value = '  -2.21%  '

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.99'
$ws.Range("E14").Value = '  -0.17%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.415.92'
$ws.Range("E15").Value = '  +0.71%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '3.764.42'
$ws.Range("E16").Value = '  +0.46%  '

# Row 17: Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.47'
$ws.Range("E17").Value = '  +2.86%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '67.632.18'
$ws.Range("E18").Value = '  -1.15%  '

# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.02'
$ws.Range("E19").Value = '  +0.66%  '

# Row 20: TRON
$ws.Range("E20").Value = '  -1.13%  '

# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  -6.61%  '

# Row 22: BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.55'
$ws.Range("E22").Value = '  -1.28%  '

# Row 23: Polygon
$ws.Range("E23").Value = '  +0.02%  '

# Row 24: PEPE
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").Value = '  +4.93%  '

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.26'
$ws.Range("E25").Value = '  -0.88%  '

# Row 26: InternetComputer(DFINITY)
$ws.Range("E26").Value = '  +0.51%  '

# Row 27: Fetch.AI
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.13'
$ws.Range("E27").Value = '  -1.98%  '

# Row 28: Dai
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.17%  '

# Row 29: RenderToken
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.01'
$ws.Range("E29").Value = '  -0.45%  '

# Row 30: PancakeSwap
$ws.Range("E30").Value = '  -0.18%  '

# Row 31: ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  +2.87%  '

# Row 32: NEARProtocol
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.21'
$ws.Range("E32").Value = '  -1.60%  '

# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.66'
$ws.Range("E33").Value = '  -0.80%  '

# Row 34: Binance-PegBSC-USD
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.03%  '

# Row 35: Aptos
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.09'
$ws.Range("E35").Value = '  -0.73%  '

# Row 36: Hedera
$ws.Range("E36").Value = '  -0.10%  '

# Row 37: dogwifhat
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.35'
$ws.Range("E37").Value = '  -0.47%  '

# Row 38: Kaspa
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("E38").Value = '  -0.08%  '

# Row 39: Mantle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.994'
$ws.Range("E39").Value = '  -0.44%  '

# Row 40: Filecoin
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("E40").Value = '  -0.37%  '

# Row 41: FirstDigitalUSD
$ws.Range("E41").Value = '  +0.01%  '

# Row 42: Arweave
$ws.Range("B42").Value = 'Arweave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.90'
$ws.Range("E42").Value = '  +5.52%  '

# Row 43: USDe
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.01%  '

# Row 44: OKB
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.27'
$ws.Range("E44").Value = '  +3.59%  '

# Row 45: TheGraph
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.298'
$ws.Range("E45").Value = '  -0.99%  '

# Row 46: Monero
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '149.77'
$ws.Range("E46").Value = '  +3.07%  '

# Row 47: Cosmos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.31'
$ws.Range("E47").Value = '  -1.90%  '

# Row 48: Bittensor
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '392.58'
$ws.Range("E48").Value = '  +0.74%  '

# Row 49: Stacks
$ws.Range("E49").Value = '  -4.70%  '

# Row 50: EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.32'
$ws.Range("E50").Value = '  +1.13%  '

# Row 51: Maker
$ws.Range("D51").Value = '2.721.47'
$ws.Range("E51").Value = '  -1.34%  '
